$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 (C35 passiv): clear all data cells except the label in A6 ---
$ws.Range("B6:K6").ClearContents()

# --- Row 7 (C35 Moderat): clear all data cells except the label in A7 ---
$ws.Range("B7:K7").ClearContents()

# --- Row 8 (C45 Aggressiv): clear all data cells except the label in A8 ---
$ws.Range("B8:K8").ClearContents()

# --- Row 9 (Facade sort glat): update data values and add date in L9 ---
$ws.Range("B9").Value = 2552.442
$ws.Range("C9").Value = 224
$ws.Range("D9").Value = 38.5
$ws.Range("E9").Value = 6.9
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 10.199999999999999
$ws.Range("H9").Value = 11.1
$ws.Range("I9").Value = 1.78
$ws.Range("J9").NumberFormat = "0.00E+00"
$ws.Range("J9").Value = 0.35499999999999998
$ws.Range("K9").Value = -9.8699999999999992
$ws.Range("L9").NumberFormat = "mmm-yy"
$ws.Range("L9").Value = 46054

# --- Row 10 (Facade sort frilagt): update data values and add date in L10 ---
$ws.Range("B10").Value = 2552.442
$ws.Range("C10").Value = 224
$ws.Range("D10").Value = 38.5
$ws.Range("E10").Value = 6.9
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 10.199999999999999
$ws.Range("H10").Value = 11.1
$ws.Range("I10").Value = 1.78
$ws.Range("J10").NumberFormat = "0.00E+00"
$ws.Range("J10").Value = 0.35499999999999998
$ws.Range("K10").Value = -9.8699999999999992
$ws.Range("L10").NumberFormat = "mmm-yy"
$ws.Range("L10").Value = 46054

# --- Row 11 (Facade grå glat): clear all data cells except the label in A11 ---
$ws.Range("B11:K11").ClearContents()

# --- Row 12 (Facade hvid glat): clear all data cells except the label in A12 ---
$ws.Range("B12:K12").ClearContents()

# --- Row 13 (Facade hvid frilagt): clear all data cells except the label in A13 ---
$ws.Range("B13:K13").ClearContents()

# --- Row 14 (Facade stennungssund frilagt): clear all data cells except the label in A14 ---
$ws.Range("B14:K14").ClearContents()

# --- Row 15: rename "Isolering" to "Rockwool isolering" and update data values ---
$ws.Range("A15").Value = "Rockwool isolering"
$ws.Range("B15").Value = 60
$ws.Range("C15").NumberFormat = "0.00E+00"
$ws.Range("C15").Value = 0.436
$ws.Range("D15").NumberFormat = "0.00E+00"
$ws.Range("D15").Value = 0.0083400000000000002
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").NumberFormat = "0.00E+00"
$ws.Range("G15").Value = 0.24
$ws.Range("H15").NumberFormat = "0.00E+00"
$ws.Range("H15").Value = 0.27300000000000002
$ws.Range("I15").NumberFormat = "0.00E+00"
$ws.Range("I15").Value = 0.041799999999999997
$ws.Range("J15").NumberFormat = "0.00E+00"
$ws.Range("J15").Value = 0.0083300000000000006
$ws.Range("K15").NumberFormat = "0.00E+00"
$ws.Range("K15").Value = -0.13600000000000001
$ws.Range("L15").NumberFormat = "mmm-yy"
$ws.Range("L15").Value = 46054

# --- Update the active cell/selection to match the saved view state ---
$ws.Range("B8").Select()
